# Season-record columns: Wins / Losses / Ties
#
# The previous scraper only pulled team/player statistics (A:AC). This adds
# the season record - Wins, Losses, Ties - as three new trailing columns
# (AD:AF), with the same constant record (88-74-0) stamped on every player
# row, matching the rest of the sheet's "one record value per row" layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new headers the same look as the rest of the header row (bold,
# centered/top aligned, boxed) by copying the format from an existing
# header cell instead of re-deriving it by hand.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows (2..46): Wins=88, Losses=74, Ties=0 --------------------------
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
